# Auto-generated edit script applying the Phoenix_Profits.xlsx diff
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(3, 8).Value = 31818.75
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 10).Value = 31818.75
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 12).Value = 31818.75
$ws.Cells.Item(3, 14).Value = -32046.75

$ws.Cells.Item(18, 8).Value = 58824790
$ws.Cells.Item(18, 9).Value = 66668028
$ws.Cells.Item(18, 10).Value = 522
$ws.Cells.Item(18, 11).Value = 66668028
$ws.Cells.Item(18, 12).Value = 522
$ws.Cells.Item(18, 13).Value = -66667744
$ws.Cells.Item(18, 14).Value = -1090

$ws.Cells.Item(19, 8).Value = 5060.143
$ws.Cells.Item(19, 9).Value = 3535.3572
$ws.Cells.Item(19, 10).Value = 8109.7144
$ws.Cells.Item(19, 11).Value = 3535.3572
$ws.Cells.Item(19, 12).Value = 8109.7144
$ws.Cells.Item(19, 13).Value = -3360.3572
$ws.Cells.Item(19, 14).Value = -8459.714400000001

$ws.Cells.Item(28, 8).Value = 1163.1765
$ws.Cells.Item(28, 9).Value = 1112.4286
$ws.Cells.Item(28, 10).Value = 1400
$ws.Cells.Item(28, 11).Value = 1112.4286
$ws.Cells.Item(28, 12).Value = 1400
$ws.Cells.Item(28, 13).Value = -627.4286
$ws.Cells.Item(28, 14).Value = -2370

$ws.Cells.Item(33, 8).Value = 15152799
$ws.Cells.Item(33, 9).Value = 22223076
$ws.Cells.Item(33, 10).Value = 2204.4285
$ws.Cells.Item(33, 11).Value = 22223076
$ws.Cells.Item(33, 12).Value = 2204.4285
$ws.Cells.Item(33, 13).Value = -22222847
$ws.Cells.Item(33, 14).Value = -2662.4285

$ws.Cells.Item(39, 8).Value = 66955.25
$ws.Cells.Item(39, 9).Value = 66955.25
$ws.Cells.Item(39, 10).Value = 0
$ws.Cells.Item(39, 11).Value = 200865.75
$ws.Cells.Item(39, 12).Value = 0
$ws.Cells.Item(39, 13).Value = -200569.75

$ws.Cells.Item(61, 8).Value = 0
$ws.Cells.Item(61, 9).Value = 0
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 11).Value = 0
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(61, 13).ClearContents()

$ws.Cells.Item(62, 8).Value = 1576.2307
$ws.Cells.Item(62, 9).Value = 1640.9166
$ws.Cells.Item(62, 10).Value = 800
$ws.Cells.Item(62, 11).Value = 1640.9166
$ws.Cells.Item(62, 12).Value = 800
$ws.Cells.Item(62, 13).Value = -1016.9166
$ws.Cells.Item(62, 14).Value = -2048

$ws.Cells.Item(65, 8).Value = 1576.2307
$ws.Cells.Item(65, 9).Value = 1640.9166
$ws.Cells.Item(65, 10).Value = 800
$ws.Cells.Item(65, 11).Value = 8204.583000000001
$ws.Cells.Item(65, 12).Value = 4000
$ws.Cells.Item(65, 13).Value = -5084.583000000001
$ws.Cells.Item(65, 14).Value = -10240

$ws.Cells.Item(102, 8).Value = 31818.75
$ws.Cells.Item(102, 9).Value = 0
$ws.Cells.Item(102, 10).Value = 31818.75
$ws.Cells.Item(102, 11).Value = 0
$ws.Cells.Item(102, 12).Value = 31818.75
$ws.Cells.Item(102, 14).Value = -38308.75

$ws.Cells.Item(107, 8).Value = 997.0323
$ws.Cells.Item(107, 9).Value = 1088.8334
$ws.Cells.Item(107, 10).Value = 682.2857
$ws.Cells.Item(107, 11).Value = 1088.8334
$ws.Cells.Item(107, 12).Value = 682.2857
$ws.Cells.Item(107, 13).Value = 831.1666
$ws.Cells.Item(107, 14).Value = -4522.2857

$ws.Cells.Item(116, 8).Value = 7191.853
$ws.Cells.Item(116, 9).Value = 6669
$ws.Cells.Item(116, 10).Value = 8644.223
$ws.Cells.Item(116, 11).Value = 6669
$ws.Cells.Item(116, 12).Value = 8644.223
$ws.Cells.Item(116, 13).Value = -3227
$ws.Cells.Item(116, 14).Value = -15528.223

$ws.Cells.Item(129, 8).Value = 2099.7778
$ws.Cells.Item(129, 9).Value = 856.8570999999999
$ws.Cells.Item(129, 10).Value = 6450
$ws.Cells.Item(129, 11).Value = 2570.5713
$ws.Cells.Item(129, 12).Value = 19350
$ws.Cells.Item(129, 13).Value = 2429.4287
$ws.Cells.Item(129, 14).Value = -29350

$ws.Cells.Item(131, 8).Value = 7364.8667
$ws.Cells.Item(131, 9).Value = 3547.7
$ws.Cells.Item(131, 10).Value = 14999.2
$ws.Cells.Item(131, 11).Value = 10643.1
$ws.Cells.Item(131, 12).Value = 44997.60000000001
$ws.Cells.Item(131, 13).Value = -5603.099999999999
$ws.Cells.Item(131, 14).Value = -55077.60000000001

$ws.Cells.Item(132, 8).Value = 2715.0417
$ws.Cells.Item(132, 9).Value = 2715.0417
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 8145.125100000001
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -5615.125100000001

$ws.Cells.Item(137, 8).Value = 1864.3572
$ws.Cells.Item(137, 9).Value = 1620.2
$ws.Cells.Item(137, 10).Value = 2474.75
$ws.Cells.Item(137, 11).Value = 4860.6
$ws.Cells.Item(137, 12).Value = 7424.25
$ws.Cells.Item(137, 13).Value = -2310.6
$ws.Cells.Item(137, 14).Value = -12524.25

$ws.Cells.Item(138, 8).Value = 3170.946
$ws.Cells.Item(138, 9).Value = 1824.4546
$ws.Cells.Item(138, 10).Value = 3740.6155
$ws.Cells.Item(138, 11).Value = 5473.3638
$ws.Cells.Item(138, 12).Value = 11221.8465
$ws.Cells.Item(138, 13).Value = -333.3638000000001
$ws.Cells.Item(138, 14).Value = -21501.8465

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1180.3448
$ws.Cells.Item(2, 9).Value = 1304.3182
$ws.Cells.Item(2, 10).Value = 790.7143
$ws.Cells.Item(2, 11).Value = 1304.3182
$ws.Cells.Item(2, 12).Value = 790.7143
$ws.Cells.Item(2, 13).Value = -1191.3182
$ws.Cells.Item(2, 14).Value = -1016.7143

$ws.Cells.Item(23, 8).Value = 18000
$ws.Cells.Item(23, 9).Value = 0
$ws.Cells.Item(23, 10).Value = 18000
$ws.Cells.Item(23, 11).Value = 0
$ws.Cells.Item(23, 12).Value = 18000
$ws.Cells.Item(23, 14).Value = -18518

$ws.Cells.Item(43, 8).Value = 20121.8
$ws.Cells.Item(43, 9).Value = 0
$ws.Cells.Item(43, 10).Value = 20121.8
$ws.Cells.Item(43, 11).Value = 0
$ws.Cells.Item(43, 12).Value = 20121.8
$ws.Cells.Item(43, 14).Value = -20747.8

$ws.Cells.Item(45, 8).Value = 2179.25
$ws.Cells.Item(45, 9).Value = 1153.5883
$ws.Cells.Item(45, 10).Value = 4670.143
$ws.Cells.Item(45, 11).Value = 1153.5883
$ws.Cells.Item(45, 12).Value = 4670.143
$ws.Cells.Item(45, 13).Value = -776.5882999999999
$ws.Cells.Item(45, 14).Value = -5424.143

$ws.Cells.Item(50, 8).Value = 1996
$ws.Cells.Item(50, 9).Value = 996.3333
$ws.Cells.Item(50, 10).Value = 2995.6667
$ws.Cells.Item(50, 11).Value = 996.3333
$ws.Cells.Item(50, 12).Value = 2995.6667
$ws.Cells.Item(50, 13).Value = -282.3333
$ws.Cells.Item(50, 14).Value = -4423.6667

$ws.Cells.Item(116, 8).Value = 1180.3448
$ws.Cells.Item(116, 9).Value = 1304.3182
$ws.Cells.Item(116, 10).Value = 790.7143
$ws.Cells.Item(116, 11).Value = 1304.3182
$ws.Cells.Item(116, 12).Value = 790.7143
$ws.Cells.Item(116, 13).Value = 989.6818000000001
$ws.Cells.Item(116, 14).Value = -5378.7143

$ws.Cells.Item(122, 8).Value = 1531.375
$ws.Cells.Item(122, 9).Value = 1274.5
$ws.Cells.Item(122, 10).Value = 1788.25
$ws.Cells.Item(122, 11).Value = 3823.5
$ws.Cells.Item(122, 12).Value = 5364.75
$ws.Cells.Item(122, 13).Value = -1373.5
$ws.Cells.Item(122, 14).Value = -10264.75

$ws.Cells.Item(132, 8).Value = 2433.0908
$ws.Cells.Item(132, 9).Value = 2326.4333
$ws.Cells.Item(132, 10).Value = 3499.6667
$ws.Cells.Item(132, 11).Value = 6979.2999
$ws.Cells.Item(132, 12).Value = 10499.0001
$ws.Cells.Item(132, 13).Value = -4449.2999
$ws.Cells.Item(132, 14).Value = -15559.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1180.3448
$ws.Cells.Item(3, 9).Value = 1304.3182
$ws.Cells.Item(3, 10).Value = 790.7143
$ws.Cells.Item(3, 11).Value = 1304.3182
$ws.Cells.Item(3, 12).Value = 790.7143
$ws.Cells.Item(3, 13).Value = -1190.3182
$ws.Cells.Item(3, 14).Value = -1018.7143

$ws.Cells.Item(20, 8).Value = 2695.3333
$ws.Cells.Item(20, 9).Value = 2391.3333
$ws.Cells.Item(20, 10).Value = 2999.3333
$ws.Cells.Item(20, 11).Value = 2391.3333
$ws.Cells.Item(20, 12).Value = 2999.3333
$ws.Cells.Item(20, 13).Value = -2144.3333
$ws.Cells.Item(20, 14).Value = -3493.3333

$ws.Cells.Item(26, 8).Value = 24482.4
$ws.Cells.Item(26, 9).Value = 24482.4
$ws.Cells.Item(26, 10).Value = 0
$ws.Cells.Item(26, 11).Value = 24482.4
$ws.Cells.Item(26, 12).Value = 0
$ws.Cells.Item(26, 13).Value = -24190.4

$ws.Cells.Item(86, 8).Value = 52633200
$ws.Cells.Item(86, 9).Value = 100001224
$ws.Cells.Item(86, 10).Value = 2056.7778
$ws.Cells.Item(86, 11).Value = 100001224
$ws.Cells.Item(86, 12).Value = 2056.7778
$ws.Cells.Item(86, 13).Value = -100000101
$ws.Cells.Item(86, 14).Value = -4302.7778

$ws.Cells.Item(89, 8).Value = 52633200
$ws.Cells.Item(89, 9).Value = 100001224
$ws.Cells.Item(89, 10).Value = 2056.7778
$ws.Cells.Item(89, 11).Value = 500006120
$ws.Cells.Item(89, 12).Value = 10283.889
$ws.Cells.Item(89, 13).Value = -500000504
$ws.Cells.Item(89, 14).Value = -21515.889

$ws.Cells.Item(94, 8).Value = 61924.934
$ws.Cells.Item(94, 9).Value = 1549.375
$ws.Cells.Item(94, 10).Value = 130925.57
$ws.Cells.Item(94, 11).Value = 1549.375
$ws.Cells.Item(94, 12).Value = 130925.57
$ws.Cells.Item(94, 13).Value = -1098.375
$ws.Cells.Item(94, 14).Value = -131827.57

$ws.Cells.Item(134, 8).Value = 5193.375
$ws.Cells.Item(134, 9).Value = 4190.4
$ws.Cells.Item(134, 10).Value = 6865
$ws.Cells.Item(134, 11).Value = 12571.2
$ws.Cells.Item(134, 12).Value = 20595
$ws.Cells.Item(134, 13).Value = -10036.2
$ws.Cells.Item(134, 14).Value = -25665

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1238.4286
$ws.Cells.Item(31, 9).Value = 944.25
$ws.Cells.Item(31, 10).Value = 2545.889
$ws.Cells.Item(31, 11).Value = 944.25
$ws.Cells.Item(31, 12).Value = 2545.889
$ws.Cells.Item(31, 13).Value = -649.25
$ws.Cells.Item(31, 14).Value = -3135.889

$ws.Cells.Item(34, 8).Value = 1238.4286
$ws.Cells.Item(34, 9).Value = 944.25
$ws.Cells.Item(34, 10).Value = 2545.889
$ws.Cells.Item(34, 11).Value = 944.25
$ws.Cells.Item(34, 12).Value = 2545.889
$ws.Cells.Item(34, 13).Value = -742.25
$ws.Cells.Item(34, 14).Value = -2949.889

$ws.Cells.Item(62, 8).Value = 74515.07000000001
$ws.Cells.Item(62, 9).Value = 173042.17
$ws.Cells.Item(62, 10).Value = 8830.333000000001
$ws.Cells.Item(62, 11).Value = 173042.17
$ws.Cells.Item(62, 12).Value = 8830.333000000001
$ws.Cells.Item(62, 13).Value = -172418.17
$ws.Cells.Item(62, 14).Value = -10078.333

$ws.Cells.Item(65, 8).Value = 74515.07000000001
$ws.Cells.Item(65, 9).Value = 173042.17
$ws.Cells.Item(65, 10).Value = 8830.333000000001
$ws.Cells.Item(65, 11).Value = 865210.8500000001
$ws.Cells.Item(65, 12).Value = 44151.665
$ws.Cells.Item(65, 13).Value = -862090.8500000001
$ws.Cells.Item(65, 14).Value = -50391.665

$ws.Cells.Item(134, 8).Value = 3176153
$ws.Cells.Item(134, 9).Value = 3368526.2
$ws.Cells.Item(134, 10).Value = 1994.5
$ws.Cells.Item(134, 11).Value = 10105578.6
$ws.Cells.Item(134, 12).Value = 5983.5
$ws.Cells.Item(134, 13).Value = -10103043.6
$ws.Cells.Item(134, 14).Value = -11053.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(8, 8).Value = 2104
$ws.Cells.Item(8, 9).Value = 2104
$ws.Cells.Item(8, 10).Value = 0
$ws.Cells.Item(8, 11).Value = 6312
$ws.Cells.Item(8, 12).Value = 0
$ws.Cells.Item(8, 13).Value = -6173

$ws.Cells.Item(37, 8).Value = 129996
$ws.Cells.Item(37, 9).Value = 0
$ws.Cells.Item(37, 10).Value = 129996
$ws.Cells.Item(37, 11).Value = 0
$ws.Cells.Item(37, 12).Value = 389988
$ws.Cells.Item(37, 14).Value = -390212

$ws.Cells.Item(54, 8).Value = 14559.857
$ws.Cells.Item(54, 9).Value = 5452
$ws.Cells.Item(54, 10).Value = 18203
$ws.Cells.Item(54, 11).Value = 16356
$ws.Cells.Item(54, 12).Value = 54609
$ws.Cells.Item(54, 13).Value = -15797
$ws.Cells.Item(54, 14).Value = -55727

$ws.Cells.Item(76, 8).Value = 6915
$ws.Cells.Item(76, 9).Value = 0
$ws.Cells.Item(76, 10).Value = 6915
$ws.Cells.Item(76, 11).Value = 0
$ws.Cells.Item(76, 12).Value = 20745
$ws.Cells.Item(76, 14).Value = -21511

$ws.Cells.Item(79, 8).Value = 6915
$ws.Cells.Item(79, 9).Value = 0
$ws.Cells.Item(79, 10).Value = 6915
$ws.Cells.Item(79, 11).Value = 0
$ws.Cells.Item(79, 12).Value = 20745
$ws.Cells.Item(79, 14).Value = -23397

$ws.Cells.Item(92, 8).Value = 275.7143
$ws.Cells.Item(92, 9).Value = 118.5
$ws.Cells.Item(92, 10).Value = 338.6
$ws.Cells.Item(92, 11).Value = 355.5
$ws.Cells.Item(92, 12).Value = 1015.8
$ws.Cells.Item(92, 13).Value = 892.5
$ws.Cells.Item(92, 14).Value = -3511.8

$ws.Cells.Item(93, 8).Value = 14374.692
$ws.Cells.Item(93, 9).Value = 0
$ws.Cells.Item(93, 10).Value = 14374.692
$ws.Cells.Item(93, 11).Value = 0
$ws.Cells.Item(93, 12).Value = 43124.076
$ws.Cells.Item(93, 14).Value = -46868.076

$ws.Cells.Item(131, 8).Value = 1913.0625
$ws.Cells.Item(131, 9).Value = 2033.1
$ws.Cells.Item(131, 10).Value = 1713
$ws.Cells.Item(131, 11).Value = 6099.299999999999
$ws.Cells.Item(131, 12).Value = 5139
$ws.Cells.Item(131, 13).Value = -1059.299999999999
$ws.Cells.Item(131, 14).Value = -15219

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 3040.3076
$ws.Cells.Item(80, 9).Value = 2356.8
$ws.Cells.Item(80, 10).Value = 3467.5
$ws.Cells.Item(80, 11).Value = 2356.8
$ws.Cells.Item(80, 12).Value = 3467.5
$ws.Cells.Item(80, 13).Value = -1358.8
$ws.Cells.Item(80, 14).Value = -5463.5

$ws.Cells.Item(83, 8).Value = 3040.3076
$ws.Cells.Item(83, 9).Value = 2356.8
$ws.Cells.Item(83, 10).Value = 3467.5
$ws.Cells.Item(83, 11).Value = 11784
$ws.Cells.Item(83, 12).Value = 17337.5
$ws.Cells.Item(83, 13).Value = -6792
$ws.Cells.Item(83, 14).Value = -27321.5

$ws.Cells.Item(113, 8).Value = 11672.23
$ws.Cells.Item(113, 9).Value = 12966
$ws.Cells.Item(113, 10).Value = 4556.5
$ws.Cells.Item(113, 11).Value = 12966
$ws.Cells.Item(113, 12).Value = 4556.5
$ws.Cells.Item(113, 13).Value = -10796
$ws.Cells.Item(113, 14).Value = -8896.5

$ws.Cells.Item(122, 8).Value = 84416.84
$ws.Cells.Item(122, 9).Value = 93777.61
$ws.Cells.Item(122, 10).Value = 12651
$ws.Cells.Item(122, 11).Value = 281332.83
$ws.Cells.Item(122, 12).Value = 37953
$ws.Cells.Item(122, 13).Value = -278882.83
$ws.Cells.Item(122, 14).Value = -42853

$ws.Cells.Item(126, 8).Value = 4742.6924
$ws.Cells.Item(126, 9).Value = 3521.625
$ws.Cells.Item(126, 10).Value = 6696.4
$ws.Cells.Item(126, 11).Value = 10564.875
$ws.Cells.Item(126, 12).Value = 20089.2
$ws.Cells.Item(126, 13).Value = -8094.875
$ws.Cells.Item(126, 14).Value = -25029.2

$ws.Cells.Item(132, 8).Value = 3800.3333
$ws.Cells.Item(132, 9).Value = 3918.0435
$ws.Cells.Item(132, 10).Value = 3123.5
$ws.Cells.Item(132, 11).Value = 11754.1305
$ws.Cells.Item(132, 12).Value = 9370.5
$ws.Cells.Item(132, 13).Value = -9224.130500000001
$ws.Cells.Item(132, 14).Value = -14430.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(20, 8).Value = 2399.8
$ws.Cells.Item(20, 9).Value = 2399.8
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 11).Value = 2399.8
$ws.Cells.Item(20, 12).Value = 0
$ws.Cells.Item(20, 13).Value = -2173.8

$ws.Cells.Item(40, 8).Value = 4012.7368
$ws.Cells.Item(40, 9).Value = 3842.6667
$ws.Cells.Item(40, 10).Value = 4304.2856
$ws.Cells.Item(40, 11).Value = 3842.6667
$ws.Cells.Item(40, 12).Value = 4304.2856
$ws.Cells.Item(40, 13).Value = -3706.6667
$ws.Cells.Item(40, 14).Value = -4576.2856

$ws.Cells.Item(46, 8).Value = 2860.524
$ws.Cells.Item(46, 9).Value = 790.875
$ws.Cells.Item(46, 10).Value = 4134.154
$ws.Cells.Item(46, 11).Value = 790.875
$ws.Cells.Item(46, 12).Value = 4134.154
$ws.Cells.Item(46, 13).Value = -602.875
$ws.Cells.Item(46, 14).Value = -4510.154

$ws.Cells.Item(61, 8).Value = 9304.895
$ws.Cells.Item(61, 9).Value = 9304.895
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 11).Value = 9304.895
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(61, 13).Value = -9102.895

$ws.Cells.Item(68, 8).Value = 3932.3845
$ws.Cells.Item(68, 9).Value = 2556.4546
$ws.Cells.Item(68, 10).Value = 11500
$ws.Cells.Item(68, 11).Value = 2556.4546
$ws.Cells.Item(68, 12).Value = 11500
$ws.Cells.Item(68, 13).Value = -1807.4546
$ws.Cells.Item(68, 14).Value = -12998

$ws.Cells.Item(71, 8).Value = 3932.3845
$ws.Cells.Item(71, 9).Value = 2556.4546
$ws.Cells.Item(71, 10).Value = 11500
$ws.Cells.Item(71, 11).Value = 12782.273
$ws.Cells.Item(71, 12).Value = 57500
$ws.Cells.Item(71, 13).Value = -9038.273000000001
$ws.Cells.Item(71, 14).Value = -64988

$ws.Cells.Item(82, 8).Value = 1796.7778
$ws.Cells.Item(82, 9).Value = 2217.889
$ws.Cells.Item(82, 10).Value = 1375.6666
$ws.Cells.Item(82, 11).Value = 2217.889
$ws.Cells.Item(82, 12).Value = 1375.6666
$ws.Cells.Item(82, 13).Value = -1856.889
$ws.Cells.Item(82, 14).Value = -2097.6666

$ws.Cells.Item(85, 8).Value = 1796.7778
$ws.Cells.Item(85, 9).Value = 2217.889
$ws.Cells.Item(85, 10).Value = 1375.6666
$ws.Cells.Item(85, 11).Value = 2217.889
$ws.Cells.Item(85, 12).Value = 1375.6666
$ws.Cells.Item(85, 13).Value = -969.8890000000001
$ws.Cells.Item(85, 14).Value = -3871.6666

$ws.Cells.Item(100, 8).Value = 4655.9443
$ws.Cells.Item(100, 9).Value = 4146
$ws.Cells.Item(100, 10).Value = 5675.8335
$ws.Cells.Item(100, 11).Value = 4146
$ws.Cells.Item(100, 12).Value = 5675.8335
$ws.Cells.Item(100, 13).Value = -3605
$ws.Cells.Item(100, 14).Value = -6757.8335

$ws.Cells.Item(113, 8).Value = 9304.895
$ws.Cells.Item(113, 9).Value = 9304.895
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 9304.895
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = -7134.895

$ws.Cells.Item(132, 8).Value = 19708.916
$ws.Cells.Item(132, 9).Value = 17409.727
$ws.Cells.Item(132, 10).Value = 45000
$ws.Cells.Item(132, 11).Value = 52229.181
$ws.Cells.Item(132, 12).Value = 135000
$ws.Cells.Item(132, 13).Value = -49699.181
$ws.Cells.Item(132, 14).Value = -140060

$ws.Cells.Item(136, 8).Value = 3225.3215
$ws.Cells.Item(136, 9).Value = 2532.6
$ws.Cells.Item(136, 10).Value = 8998
$ws.Cells.Item(136, 11).Value = 7597.799999999999
$ws.Cells.Item(136, 12).Value = 26994
$ws.Cells.Item(136, 13).Value = -5047.799999999999
$ws.Cells.Item(136, 14).Value = -32094

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(52, 8).Value = 34994.7
$ws.Cells.Item(52, 9).Value = 20064.75
$ws.Cells.Item(52, 10).Value = 44948
$ws.Cells.Item(52, 11).Value = 20064.75
$ws.Cells.Item(52, 12).Value = 44948
$ws.Cells.Item(52, 13).Value = -19838.75
$ws.Cells.Item(52, 14).Value = -45400

$ws.Cells.Item(81, 8).Value = 11410310
$ws.Cells.Item(81, 9).Value = 16596034
$ws.Cells.Item(81, 10).Value = 1719
$ws.Cells.Item(81, 11).Value = 33192068
$ws.Cells.Item(81, 12).Value = 3438
$ws.Cells.Item(81, 13).Value = -33191007
$ws.Cells.Item(81, 14).Value = -5560

$ws.Cells.Item(84, 8).Value = 11410310
$ws.Cells.Item(84, 9).Value = 16596034
$ws.Cells.Item(84, 10).Value = 1719
$ws.Cells.Item(84, 11).Value = 165960340
$ws.Cells.Item(84, 12).Value = 17190
$ws.Cells.Item(84, 13).Value = -165955036
$ws.Cells.Item(84, 14).Value = -27798

$ws.Cells.Item(96, 8).Value = 5686.5557
$ws.Cells.Item(96, 9).Value = 10866
$ws.Cells.Item(96, 10).Value = 3096.8333
$ws.Cells.Item(96, 11).Value = 10866
$ws.Cells.Item(96, 12).Value = 3096.8333
$ws.Cells.Item(96, 13).Value = -9493
$ws.Cells.Item(96, 14).Value = -5842.8333

$ws.Cells.Item(100, 8).Value = 175.5
$ws.Cells.Item(100, 9).Value = 166.5
$ws.Cells.Item(100, 10).Value = 193.5
$ws.Cells.Item(100, 11).Value = 333
$ws.Cells.Item(100, 12).Value = 387
$ws.Cells.Item(100, 13).Value = 208
$ws.Cells.Item(100, 14).Value = -1469

$ws.Cells.Item(120, 8).Value = 69300
$ws.Cells.Item(120, 9).Value = 0
$ws.Cells.Item(120, 10).Value = 69300
$ws.Cells.Item(120, 11).Value = 0
$ws.Cells.Item(120, 12).Value = 69300
$ws.Cells.Item(120, 14).Value = -78976

$ws.Cells.Item(126, 8).Value = 3374
$ws.Cells.Item(126, 9).Value = 3332.3333
$ws.Cells.Item(126, 10).Value = 3499
$ws.Cells.Item(126, 11).Value = 9996.999899999999
$ws.Cells.Item(126, 12).Value = 10497
$ws.Cells.Item(126, 13).Value = -7526.999899999999
$ws.Cells.Item(126, 14).Value = -15437

$ws.Cells.Item(132, 8).Value = 2408.0488
$ws.Cells.Item(132, 9).Value = 2395.3948
$ws.Cells.Item(132, 10).Value = 2568.3333
$ws.Cells.Item(132, 11).Value = 7186.1844
$ws.Cells.Item(132, 12).Value = 7704.999899999999
$ws.Cells.Item(132, 13).Value = -4656.1844
$ws.Cells.Item(132, 14).Value = -12764.9999

$ws.Cells.Item(136, 8).Value = 1032.8667
$ws.Cells.Item(136, 9).Value = 719.28
$ws.Cells.Item(136, 10).Value = 2600.8
$ws.Cells.Item(136, 11).Value = 2157.84
$ws.Cells.Item(136, 12).Value = 7802.400000000001
$ws.Cells.Item(136, 13).Value = 392.1599999999999
$ws.Cells.Item(136, 14).Value = -12902.4
